# Insert a brand-new weekly price row for "Acelga" (Macroferia Regional de
# Talca) at row 385, pushing the previously-existing rows 385:522 down to
# 386:523 (dimension grows from A1:R522 to A1:R523).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 385 and below down by one row.
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new record.
$ws.Range("A385").Value = 5
$ws.Range("B385").Value = "Macroferia Regional de Talca"
$ws.Range("C385").Value = "Maule"
$ws.Range("D385").Value = 45229
$ws.Range("E385").Value = 7
$ws.Range("F385").Value = 100112009
$ws.Range("G385").Value = "Acelga"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 300
$ws.Range("K385").Value = 2000
$ws.Range("L385").Value = 2000
$ws.Range("M385").Value = 2000
$ws.Range("N385").Value = "$/docena de atados (4 kilos)"
$ws.Range("O385").Value = "Región del Maule"
$ws.Range("P385").Value = 500
$ws.Range("Q385").Value = 4
$ws.Range("R385").Value = "Hortaliza"
